$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 92

$rng = $ws.Range("A" + $row + ":E" + $row)
$rng.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-12-17"
$ws.Cells.Item($row, 2).Value = "Pick 4"
$ws.Cells.Item($row, 3).Value = "251217"
$ws.Cells.Item($row, 4).Value = "2-6-0-5"
$ws.Cells.Item($row, 5).Value = "2025-12-17T21:44:54.387+04:00"

$rng.Style = "Normal"
